$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 344, shifting existing rows 344-414 down to 345-415
$ws.Rows("344").Insert()

# Populate the newly inserted row 344 with its data
$ws.Range("A344").Value = 6
$ws.Range("B344").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C344").Value = "Metropolitana"
$ws.Range("D344").Value = 44511
$ws.Range("E344").Value = 13
$ws.Range("F344").Value = 100112003
$ws.Range("G344").Value = "Ajo"
$ws.Range("H344").Value = "Chino"
$ws.Range("I344").Value = "Primera"
$ws.Range("J344").Value = 2780
$ws.Range("K344").Value = 16500
$ws.Range("L344").Value = 17000
$ws.Range("M344").Value = 16770
$ws.Range("N344").Value = "$/caja 10 kilos"
$ws.Range("O344").Value = "China"
$ws.Range("P344").Value = 1677
$ws.Range("Q344").Value = 10
$ws.Range("R344").Value = "Hortaliza"
